$wb = $excel.ActiveWorkbook

# The last sheet ("Namrata 213") is the template for the new PO sheet.
$lastIndex = $wb.Worksheets.Count
$template = $wb.Worksheets.Item($lastIndex)

# Copy it to the end of the workbook -> becomes the new active/last sheet.
$template.Copy($null, $template)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Aditya 214"

# Update header row (column B header -> "Model Number")
$newSheet.Range("B1").Value = "Model Number"

# Update the single line-item row with the new product data.
$newSheet.Range("B2").Value = "DGS-1210-10P"
$newSheet.Range("C2").Value = "Dlink 8-Ports 10/100/1000Mbps POE + 2-Ports SFP 100/1000Mbps Smart Managed Switch, 65Watts "
$newSheet.Range("D2").Value = 5
$newSheet.Range("E2").Value = 9600

# Row heights: header row becomes shorter on both the template and the new sheet.
$template.Rows.Item(1).RowHeight = 28.8
$newSheet.Rows.Item(1).RowHeight = 28.8

# Restore selections as left by the author after the edit.
$template.Range("A1:F5").Select()
$newSheet.Range("G8").Select()

$cassun212 = $wb.Worksheets.Item("Cassun Electricals 212")
$cassun212.Range("I27").Select()

$newSheet.Activate()
